# TS3-Zobrazenie skladu.xlsx - pridavam vyslednu technicku dokumentaciu
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the (empty) "Dodatocne informacie" column E - this shifts
#    the old "Ocakavany vysledok" column (F) left into the E position.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).Delete()

# ---------------------------------------------------------------------
# 2) Header row 5 label changes
# ---------------------------------------------------------------------
$ws.Range("B5").Value2 = "Krok:"
$ws.Range("E5").Value2 = "Očakávaný výsledok:"

# ---------------------------------------------------------------------
# 3) Merge B3:C3 and combine the "Predpoklady:" label with its answer
#    (previously split across B3 label / C3 value).
# ---------------------------------------------------------------------
$ws.Range("B3:C3").UnMerge() | Out-Null
$ws.Range("B3").Value2 = "Predpoklady: Prístup k aplikácii s právami používateľa."
$prefixLen = ("Predpoklady: ").Length
$totalLen = $ws.Range("B3").Value2.Length
$ws.Range("B3").Characters($prefixLen + 1, $totalLen - $prefixLen).Font.Bold = $false
$ws.Range("C3").Value2 = ""
$ws.Range("B3:C3").Merge() | Out-Null

# ---------------------------------------------------------------------
# 4) Reword the palette-info cell (now E13) - group related fields onto
#    shared lines.
# ---------------------------------------------------------------------
$nl = [char]10
$ws.Range("E13").Value2 = "Systém vypíše informácie o palete:" + $nl + "Typ palety, Hmotnosť, Nadrozmernosť" + $nl + "Poškodenosť, Ponámka" + $nl + "Meno skladníka ktorý zaskladňoval paletu" + $nl + "Zákazník od ktorého palte prišla" + $nl + "Tabuľku s materiálmi na palete a ich počtami"

Write-Host "content done"

# ---------------------------------------------------------------------
# 5) Column widths (characters) - col A stays close to original, B/C/D/E
#    get re-tuned slightly (artifact of the content/formatting refresh).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 1.5625
$ws.Columns.Item(2).ColumnWidth = 6.15234375
$ws.Columns.Item(3).ColumnWidth = 40.234375
$ws.Columns.Item(4).ColumnWidth = 24.8046875
$ws.Columns.Item(5).ColumnWidth = 53.7109375

# ---------------------------------------------------------------------
# 6) Row heights (points)
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 10.15
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 13.9
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 35.25
$ws.Rows.Item(7).RowHeight = 31.5
$ws.Rows.Item(8).RowHeight = 42
$ws.Rows.Item(9).RowHeight = 44.45
$ws.Rows.Item(10).RowHeight = 20.25
$ws.Rows.Item(11).RowHeight = 69
$ws.Rows.Item(12).RowHeight = 63.75
$ws.Rows.Item(13).RowHeight = 102.75
$ws.Rows.Item(14).RowHeight = 14.45
$ws.Rows.Item(15).RowHeight = 34.15
$ws.Rows.Item(16).RowHeight = 37.15
$ws.Rows.Item(17).RowHeight = 53.45
$ws.Rows.Item(18).RowHeight = 72.6
$ws.Rows.Item(19).RowHeight = 71.45
$ws.Rows.Item(20).RowHeight = 45

# ---------------------------------------------------------------------
# 7) Page setup: landscape orientation
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 2

# ---------------------------------------------------------------------
# 8) View: drop the frozen/scrolled "topLeftCell" and update selection
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G12").Select()

Write-Host "formatting done"
